# "Excel Datei Reiter Linien ergänzt"
# Insert a new worksheet "Linie" between "QS" and "Tabelle1", give it some
# sample values, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

$qs = $wb.Worksheets.Item("QS")

# Add the new sheet right after "QS" (i.e. before "Tabelle1").
$linie = $wb.Worksheets.Add($null, $qs)
$linie.Name = "Linie"

# Sample data for the new sheet.
$linie.Range("A1").Value = 1
$linie.Range("B1").Value = 1
$linie.Range("C1").Value = 2

# Make "Linie" the active sheet/tab with the same selection as in the
# target workbook.
$linie.Activate()
$linie.Range("C12").Select()
